$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '66.151.68'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +7.16%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.016.68'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +4.39%  '

$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '585.03'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.45%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.87'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +9.06%  '

$ws.Range("E7").Value = '  -0.06%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '3.011.42'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +4.25%  '

$ws.Range("E9").Value = '  +3.41%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.97'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.46%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.155'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +6.66%  '

$ws.Range("E12").Value = '  +5.60%  '

$ws.Range("E13").Value = '  +8.66%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '34.61'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +9.31%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '66.172.17'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +7.27%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '3.517.37'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +4.39%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.95'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.89%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '3.016.65'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +4.02%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '462.53'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +7.67%  '

$ws.Range("E21").Value = '  +6.74%  '

$ws.Range("E22").Value = '  +5.37%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '7.36'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +8.47%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '82.16'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.09%  '

$ws.Range("E25").Value = '  +13.18%  '

$ws.Range("E26").Value = '  +5.02%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '10.70'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +7.92%  '

$ws.Range("E28").Value = '  -0.07%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.04'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +14.62%  '

$ws.Range("E30").Value = '  +15.95%  '

$ws.Range("E31").Value = '  -1.49%  '

$ws.Range("E32").Value = '  +5.33%  '

$ws.Range("E33").Value = '  +4.69%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '27.01'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +6.24%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.999'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.08%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.995'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +3.99%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '5.81'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +8.49%  '

$ws.Range("E38").Value = '  +13.07%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.06'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +9.35%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '49.48'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.37%  '

$ws.Range("E41").Value = '  +8.33%  '

$ws.Range("E42").Value = '  +13.70%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '43.55'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +10.51%  '

$ws.Range("E44").Value = '  +3.55%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '392.19'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +14.31%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.800.85'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +4.56%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0355'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.91%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '133.83'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +1.28%  '

$ws.Range("E49").Value = '  -0.02%  '

$ws.Range("E50").Value = '  +9.61%  '

$ws.Range("E51").Value = '  +4.22%  '
